$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 112994
$ws.Range("J3").Value = 112994
$ws.Range("L3").Value = 112994
$ws.Range("N3").Value = -113222
$ws.Range("H9").Value = 7001140
$ws.Range("I9").Value = 479.6
$ws.Range("K9").Value = 479.6
$ws.Range("M9").Value = -310.6
$ws.Range("H40").Value = 100391.29
$ws.Range("I40").Value = 502071.66
$ws.Range("J40").Value = 3988
$ws.Range("K40").Value = 502071.66
$ws.Range("L40").Value = 3988
$ws.Range("M40").Value = -501896.66
$ws.Range("N40").Value = -4338
$ws.Range("H43").Value = 19277.445
$ws.Range("I43").Value = 17999.666
$ws.Range("J43").Value = 19916.334
$ws.Range("K43").Value = 17999.666
$ws.Range("L43").Value = 19916.334
$ws.Range("M43").Value = -17930.666
$ws.Range("N43").Value = -20054.334
$ws.Range("H62").Value = 4076
$ws.Range("J62").Value = 5300.0835
$ws.Range("L62").Value = 5300.0835
$ws.Range("N62").Value = -6548.0835
$ws.Range("H65").Value = 4076
$ws.Range("J65").Value = 5300.0835
$ws.Range("L65").Value = 26500.4175
$ws.Range("N65").Value = -32740.4175
$ws.Range("H86").Value = 230775170
$ws.Range("I86").Value = 166673100
$ws.Range("K86").Value = 166673100
$ws.Range("M86").Value = -166671977
$ws.Range("H87").Value = 101865.91
$ws.Range("J87").Value = 95553.8
$ws.Range("L87").Value = 95553.8
$ws.Range("N87").Value = -98049.8
$ws.Range("H89").Value = 230775170
$ws.Range("I89").Value = 166673100
$ws.Range("K89").Value = 833365500
$ws.Range("M89").Value = -833359884
$ws.Range("H90").Value = 101865.91
$ws.Range("J90").Value = 95553.8
$ws.Range("L90").Value = 286661.4
$ws.Range("N90").Value = -299141.4
$ws.Range("H92").Value = 41667530
$ws.Range("I92").Value = 45455440
$ws.Range("K92").Value = 45455440
$ws.Range("M92").Value = -45454192
$ws.Range("H95").Value = 71869.664
$ws.Range("J95").Value = 71869.664
$ws.Range("L95").Value = 71869.664
$ws.Range("N95").Value = -77361.664
$ws.Range("H99").Value = 788
$ws.Range("J99").Value = 299
$ws.Range("L99").Value = 897
$ws.Range("N99").Value = -3893
$ws.Range("I100").Value = 2166.0476
$ws.Range("J100").Value = 4970.5293
$ws.Range("K100").Value = 2166.0476
$ws.Range("L100").Value = 4970.5293
$ws.Range("M100").Value = -1625.0476
$ws.Range("N100").Value = -6052.5293
$ws.Range("H102").Value = 112994
$ws.Range("J102").Value = 112994
$ws.Range("L102").Value = 112994
$ws.Range("N102").Value = -119484
$ws.Range("H125").Value = 2768.4614
$ws.Range("I125").Value = 2284.8572
$ws.Range("J125").Value = 3332.6667
$ws.Range("K125").Value = 20563.7148
$ws.Range("L125").Value = 29994.0003
$ws.Range("M125").Value = -18103.7148
$ws.Range("N125").Value = -34914.0003
$ws.Range("H129").Value = 2523.875
$ws.Range("I129").Value = 1049.75
$ws.Range("K129").Value = 3149.25
$ws.Range("M129").Value = 1850.75
$ws.Range("H132").Value = 2725.7812
$ws.Range("I132").Value = 2145.6897
$ws.Range("K132").Value = 6437.0691
$ws.Range("M132").Value = -3907.0691
$ws.Range("H137").Value = 1881317.9
$ws.Range("I137").Value = 103873.375
$ws.Range("J137").Value = 2527661.2
$ws.Range("K137").Value = 311620.125
$ws.Range("L137").Value = 7582983.600000001
$ws.Range("M137").Value = -309070.125
$ws.Range("N137").Value = -7588083.600000001
$ws.Range("H138").Value = 4415.6
$ws.Range("I138").Value = 3632.625
$ws.Range("J138").Value = 4483.6846
$ws.Range("K138").Value = 10897.875
$ws.Range("L138").Value = 13451.0538
$ws.Range("M138").Value = -5757.875
$ws.Range("N138").Value = -23731.0538
$ws.Range("H141").Value = 4975.6
$ws.Range("I141").Value = 4781.9565
$ws.Range("J141").Value = 7202.5
$ws.Range("K141").Value = 14345.8695
$ws.Range("L141").Value = 21607.5
$ws.Range("M141").Value = -9165.869500000001
$ws.Range("N141").Value = -31967.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1905.3529
$ws.Range("I2").Value = 2049.4167
$ws.Range("J2").Value = 1559.6
$ws.Range("K2").Value = 2049.4167
$ws.Range("L2").Value = 1559.6
$ws.Range("M2").Value = -1936.4167
$ws.Range("N2").Value = -1785.6
$ws.Range("H31").Value = 54504.168
$ws.Range("J31").Value = 80756.75
$ws.Range("L31").Value = 80756.75
$ws.Range("N31").Value = -81344.75
$ws.Range("H32").Value = 23455060
$ws.Range("I32").Value = 22587394
$ws.Range("J32").Value = 62500000
$ws.Range("K32").Value = 22587394
$ws.Range("L32").Value = 62500000
$ws.Range("M32").Value = -22587107
$ws.Range("N32").Value = -62500574
$ws.Range("H74").Value = 3285.3572
$ws.Range("I74").Value = 3599.375
$ws.Range("K74").Value = 3599.375
$ws.Range("M74").Value = -2725.375
$ws.Range("H77").Value = 3285.3572
$ws.Range("I77").Value = 3599.375
$ws.Range("K77").Value = 17996.875
$ws.Range("M77").Value = -13628.875
$ws.Range("H116").Value = 1905.3529
$ws.Range("I116").Value = 2049.4167
$ws.Range("J116").Value = 1559.6
$ws.Range("K116").Value = 2049.4167
$ws.Range("L116").Value = 1559.6
$ws.Range("M116").Value = 244.5832999999998
$ws.Range("N116").Value = -6147.6
$ws.Range("H122").Value = 4010.325
$ws.Range("I122").Value = 4047.7778
$ws.Range("J122").Value = 3673.25
$ws.Range("K122").Value = 12143.3334
$ws.Range("L122").Value = 11019.75
$ws.Range("M122").Value = -9693.3334
$ws.Range("N122").Value = -15919.75
$ws.Range("H128").Value = 89995
$ws.Range("J128").Value = 89995
$ws.Range("L128").Value = 89995
$ws.Range("N128").Value = -99955
$ws.Range("H132").Value = 273590.38
$ws.Range("I132").Value = 325982.9
$ws.Range("J132").Value = 2895.6667
$ws.Range("K132").Value = 977948.7000000001
$ws.Range("L132").Value = 8687.000100000001
$ws.Range("M132").Value = -975418.7000000001
$ws.Range("N132").Value = -13747.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1905.3529
$ws.Range("I3").Value = 2049.4167
$ws.Range("J3").Value = 1559.6
$ws.Range("K3").Value = 2049.4167
$ws.Range("L3").Value = 1559.6
$ws.Range("M3").Value = -1935.4167
$ws.Range("N3").Value = -1787.6
$ws.Range("H20").Value = 1957.9656
$ws.Range("I20").Value = 1852
$ws.Range("J20").Value = 2291
$ws.Range("K20").Value = 1852
$ws.Range("L20").Value = 2291
$ws.Range("M20").Value = -1605
$ws.Range("N20").Value = -2785
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H38").Value = 84435
$ws.Range("J38").Value = 84435
$ws.Range("L38").Value = 84435
$ws.Range("N38").Value = -85267
$ws.Range("H64").Value = 941.75
$ws.Range("J64").Value = 964.1667
$ws.Range("L64").Value = 964.1667
$ws.Range("N64").Value = -1414.1667
$ws.Range("H67").Value = 941.75
$ws.Range("J67").Value = 964.1667
$ws.Range("L67").Value = 964.1667
$ws.Range("N67").Value = -2524.1667
$ws.Range("H86").Value = 3976.2
$ws.Range("I86").Value = 3415.7778
$ws.Range("J86").Value = 4816.8335
$ws.Range("K86").Value = 3415.7778
$ws.Range("L86").Value = 4816.8335
$ws.Range("M86").Value = -2292.7778
$ws.Range("N86").Value = -7062.8335
$ws.Range("H89").Value = 3976.2
$ws.Range("I89").Value = 3415.7778
$ws.Range("J89").Value = 4816.8335
$ws.Range("K89").Value = 17078.889
$ws.Range("L89").Value = 24084.1675
$ws.Range("M89").Value = -11462.889
$ws.Range("N89").Value = -35316.1675
$ws.Range("H99").Value = 3009.6667
$ws.Range("J99").Value = 3747.25
$ws.Range("L99").Value = 3747.25
$ws.Range("N99").Value = -6743.25
$ws.Range("H132").Value = 49335
$ws.Range("J132").Value = 49335
$ws.Range("L132").Value = 49335
$ws.Range("N132").Value = -59455
$ws.Range("H134").Value = 3178366.5
$ws.Range("I134").Value = 3178366.5
$ws.Range("K134").Value = 9535099.5
$ws.Range("M134").Value = -9532564.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 540
$ws.Range("J16").Value = 530
$ws.Range("L16").Value = 530
$ws.Range("N16").Value = -1104
$ws.Range("H22").Value = 10890.909
$ws.Range("I22").Value = 1343
$ws.Range("J22").Value = 27599.75
$ws.Range("K22").Value = 1343
$ws.Range("L22").Value = 27599.75
$ws.Range("M22").Value = -993
$ws.Range("N22").Value = -28299.75
$ws.Range("H31").Value = 2998.17
$ws.Range("I31").Value = 1904.5385
$ws.Range("J31").Value = 3353.6
$ws.Range("K31").Value = 1904.5385
$ws.Range("L31").Value = 3353.6
$ws.Range("M31").Value = -1609.5385
$ws.Range("N31").Value = -3943.6
$ws.Range("H34").Value = 2998.17
$ws.Range("I34").Value = 1904.5385
$ws.Range("J34").Value = 3353.6
$ws.Range("K34").Value = 1904.5385
$ws.Range("L34").Value = 3353.6
$ws.Range("M34").Value = -1702.5385
$ws.Range("N34").Value = -3757.6
$ws.Range("H58").Value = 3714.6287
$ws.Range("I58").Value = 3060.889
$ws.Range("J58").Value = 4406.8237
$ws.Range("K58").Value = 3060.889
$ws.Range("L58").Value = 4406.8237
$ws.Range("M58").Value = -2857.889
$ws.Range("N58").Value = -4812.8237
$ws.Range("H99").Value = 2601.0625
$ws.Range("I99").Value = 2496.3333
$ws.Range("J99").Value = 2735.7144
$ws.Range("K99").Value = 2496.3333
$ws.Range("L99").Value = 2735.7144
$ws.Range("M99").Value = -998.3332999999998
$ws.Range("N99").Value = -5731.7144
$ws.Range("H105").Value = 4836
$ws.Range("J105").Value = 4066.6667
$ws.Range("L105").Value = 4066.6667
$ws.Range("N105").Value = -7560.6667
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = None
$ws.Range("N109").ClearContents()
$ws.Range("H113").Value = 540
$ws.Range("J113").Value = 530
$ws.Range("L113").Value = 530
$ws.Range("N113").Value = -4870
$ws.Range("H126").Value = 2601.0625
$ws.Range("I126").Value = 2496.3333
$ws.Range("J126").Value = 2735.7144
$ws.Range("K126").Value = 7488.999899999999
$ws.Range("L126").Value = 8207.143199999999
$ws.Range("M126").Value = -5018.999899999999
$ws.Range("N126").Value = -13147.1432
$ws.Range("H132").Value = 3883
$ws.Range("I132").Value = 3740.6924
$ws.Range("J132").Value = 4279.4287
$ws.Range("K132").Value = 11222.0772
$ws.Range("L132").Value = 12838.2861
$ws.Range("M132").Value = -8692.0772
$ws.Range("N132").Value = -17898.2861
$ws.Range("H134").Value = 2182
$ws.Range("I134").Value = 2005.3182
$ws.Range("K134").Value = 6015.9546
$ws.Range("M134").Value = -3480.9546
$ws.Range("H136").Value = 3714.6287
$ws.Range("I136").Value = 3060.889
$ws.Range("J136").Value = 4406.8237
$ws.Range("K136").Value = 9182.667000000001
$ws.Range("L136").Value = 13220.4711
$ws.Range("M136").Value = -6632.667000000001
$ws.Range("N136").Value = -18320.4711

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 4329052
$ws.Range("I2").Value = 5050526
$ws.Range("K2").Value = 30303156
$ws.Range("M2").Value = -30303043
$ws.Range("H4").Value = 142308960
$ws.Range("I4").Value = 105568470
$ws.Range("J4").Value = 399492400
$ws.Range("K4").Value = 316705410
$ws.Range("L4").Value = 1198477200
$ws.Range("M4").Value = -316705298
$ws.Range("N4").Value = -1198477424
$ws.Range("H5").Value = 539.5
$ws.Range("I5").Value = 553
$ws.Range("K5").Value = 1659
$ws.Range("M5").Value = -1547
$ws.Range("H12").Value = 500056.9
$ws.Range("J12").Value = 714352.9
$ws.Range("L12").Value = 2143058.7
$ws.Range("N12").Value = -2143404.7
$ws.Range("H14").Value = 1597.2
$ws.Range("I14").Value = 1597.2
$ws.Range("K14").Value = 4791.6
$ws.Range("M14").Value = -4618.6
$ws.Range("H38").Value = 834
$ws.Range("I38").Value = 178.33333
$ws.Range("J38").Value = 1325.75
$ws.Range("K38").Value = 534.99999
$ws.Range("L38").Value = 3977.25
$ws.Range("M38").Value = -187.99999
$ws.Range("N38").Value = -4671.25
$ws.Range("H56").Value = 5978.577
$ws.Range("I56").Value = 5978.577
$ws.Range("K56").Value = 5978.577
$ws.Range("M56").Value = -5448.577
$ws.Range("H68").Value = 1796.3
$ws.Range("I68").Value = 1141.6
$ws.Range("J68").Value = 2451
$ws.Range("K68").Value = 3424.8
$ws.Range("L68").Value = 7353
$ws.Range("M68").Value = -2613.8
$ws.Range("N68").Value = -8975
$ws.Range("H71").Value = 1796.3
$ws.Range("I71").Value = 1141.6
$ws.Range("J71").Value = 2451
$ws.Range("K71").Value = 10274.4
$ws.Range("L71").Value = 22059
$ws.Range("M71").Value = -6218.4
$ws.Range("N71").Value = -30171
$ws.Range("H98").Value = 517.8
$ws.Range("J98").Value = 130.5
$ws.Range("L98").Value = 391.5
$ws.Range("N98").Value = -3387.5
$ws.Range("H114").Value = 6248.5
$ws.Range("J114").Value = 6983.3335
$ws.Range("L114").Value = 20950.0005
$ws.Range("N114").Value = -27458.0005
$ws.Range("H131").Value = 1624.6608
$ws.Range("J131").Value = 1799.2954
$ws.Range("L131").Value = 5397.8862
$ws.Range("N131").Value = -15477.8862
$ws.Range("H135").Value = 539.5
$ws.Range("I135").Value = 553
$ws.Range("K135").Value = 4977
$ws.Range("M135").Value = -2442
$ws.Range("H136").Value = 975
$ws.Range("J141").Value = 5000
$ws.Range("L141").Value = 15000
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3865.1
$ws.Range("I80").Value = 3466.6667
$ws.Range("J80").Value = 4035.8572
$ws.Range("K80").Value = 3466.6667
$ws.Range("L80").Value = 4035.8572
$ws.Range("M80").Value = -2468.6667
$ws.Range("N80").Value = -6031.8572
$ws.Range("H83").Value = 3865.1
$ws.Range("I83").Value = 3466.6667
$ws.Range("J83").Value = 4035.8572
$ws.Range("K83").Value = 17333.3335
$ws.Range("L83").Value = 20179.286
$ws.Range("M83").Value = -12341.3335
$ws.Range("N83").Value = -30163.286
$ws.Range("H97").Value = 1818.6818
$ws.Range("I97").Value = 980
$ws.Range("J97").Value = 3615.8572
$ws.Range("K97").Value = 980
$ws.Range("L97").Value = 3615.8572
$ws.Range("M97").Value = -484
$ws.Range("N97").Value = -4607.8572
$ws.Range("H102").Value = 2966.7334
$ws.Range("I102").Value = 2736.28
$ws.Range("J102").Value = 4119
$ws.Range("K102").Value = 2736.28
$ws.Range("L102").Value = 4119
$ws.Range("M102").Value = -1114.28
$ws.Range("N102").Value = -7363
$ws.Range("H113").Value = 35655.223
$ws.Range("I113").Value = 4721.5454
$ws.Range("K113").Value = 4721.5454
$ws.Range("M113").Value = -2551.5454
$ws.Range("H122").Value = 2418.6875
$ws.Range("I122").Value = 2438.75
$ws.Range("J122").Value = 2398.625
$ws.Range("K122").Value = 7316.25
$ws.Range("L122").Value = 7195.875
$ws.Range("M122").Value = -4866.25
$ws.Range("N122").Value = -12095.875
$ws.Range("H132").Value = 3344.8064
$ws.Range("I132").Value = 2985.0527
$ws.Range("K132").Value = 8955.158100000001
$ws.Range("M132").Value = -6425.158100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1387.1578
$ws.Range("I16").Value = 1257.1333
$ws.Range("J16").Value = 1874.75
$ws.Range("K16").Value = 1257.1333
$ws.Range("L16").Value = 1874.75
$ws.Range("M16").Value = -1087.1333
$ws.Range("N16").Value = -2214.75
$ws.Range("H22").Value = 2314.8572
$ws.Range("I22").Value = 1920.6
$ws.Range("J22").Value = 2533.889
$ws.Range("K22").Value = 1920.6
$ws.Range("L22").Value = 2533.889
$ws.Range("M22").Value = -1625.6
$ws.Range("N22").Value = -3123.889
$ws.Range("H27").Value = 2314.8572
$ws.Range("I27").Value = 1920.6
$ws.Range("J27").Value = 2533.889
$ws.Range("K27").Value = 1920.6
$ws.Range("L27").Value = 2533.889
$ws.Range("M27").Value = -1813.6
$ws.Range("N27").Value = -2747.889
$ws.Range("H40").Value = 3173.4
$ws.Range("I40").Value = 2655.9546
$ws.Range("J40").Value = 6968
$ws.Range("K40").Value = 2655.9546
$ws.Range("L40").Value = 6968
$ws.Range("M40").Value = -2519.9546
$ws.Range("N40").Value = -7240
$ws.Range("H46").Value = 7329.091
$ws.Range("I46").Value = 1749
$ws.Range("J46").Value = 8210.157999999999
$ws.Range("K46").Value = 1749
$ws.Range("L46").Value = 8210.157999999999
$ws.Range("M46").Value = -1561
$ws.Range("N46").Value = -8586.157999999999
$ws.Range("H61").Value = 3869.2222
$ws.Range("I61").Value = 1619.7333
$ws.Range("K61").Value = 1619.7333
$ws.Range("M61").Value = -1417.7333
$ws.Range("H113").Value = 3869.2222
$ws.Range("I113").Value = 1619.7333
$ws.Range("K113").Value = 1619.7333
$ws.Range("M113").Value = 550.2666999999999
$ws.Range("H122").Value = 28333.666
$ws.Range("I122").Value = 29000.6
$ws.Range("K122").Value = 87001.79999999999
$ws.Range("M122").Value = -84551.79999999999
$ws.Range("H132").Value = 420482.88
$ws.Range("I132").Value = 503379.6
$ws.Range("K132").Value = 1510138.8
$ws.Range("M132").Value = -1507608.8
$ws.Range("H136").Value = 7908.4165
$ws.Range("I136").Value = 9260.429
$ws.Range("J136").Value = 6015.6
$ws.Range("K136").Value = 27781.287
$ws.Range("L136").Value = 18046.8
$ws.Range("M136").Value = -25231.287
$ws.Range("N136").Value = -23146.8
$ws.Range("H141").Value = 543333
$ws.Range("J141").Value = 543333
$ws.Range("L141").Value = 543333
$ws.Range("N141").Value = -553693

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 44999.5
$ws.Range("I43").Value = 59999
$ws.Range("K43").Value = 59999
$ws.Range("M43").Value = -59850
$ws.Range("H45").Value = 20030.75
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H62").Value = 7150.25
$ws.Range("I62").Value = 8834
$ws.Range("J62").Value = 6140
$ws.Range("K62").Value = 8834
$ws.Range("L62").Value = 6140
$ws.Range("M62").Value = -8210
$ws.Range("N62").Value = -7388
$ws.Range("H65").Value = 7150.25
$ws.Range("I65").Value = 8834
$ws.Range("J65").Value = 6140
$ws.Range("K65").Value = 44170
$ws.Range("L65").Value = 30700
$ws.Range("M65").Value = -41050
$ws.Range("N65").Value = -36940
$ws.Range("H100").Value = 2542.4
$ws.Range("J100").Value = 799.6667
$ws.Range("L100").Value = 1599.3334
$ws.Range("N100").Value = -2681.3334
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = None
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H119").Value = 109994.5
$ws.Range("J119").Value = 109994.5
$ws.Range("L119").Value = 109994.5
$ws.Range("N119").Value = -119670.5
$ws.Range("H122").Value = 4015.4827
$ws.Range("I122").Value = 4279.5415
$ws.Range("J122").Value = 2748
$ws.Range("K122").Value = 12838.6245
$ws.Range("L122").Value = 8244
$ws.Range("M122").Value = -10388.6245
$ws.Range("N122").Value = -13144
$ws.Range("H124").Value = 53332.668
$ws.Range("J124").Value = 53332.668
$ws.Range("L124").Value = 53332.668
$ws.Range("N124").Value = -63152.668
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = None
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 2900.3
$ws.Range("I126").Value = 2900.3
$ws.Range("K126").Value = 8700.900000000001
$ws.Range("M126").Value = -6230.900000000001
$ws.Range("H132").Value = 30807.314
$ws.Range("I132").Value = 34533.773
$ws.Range("J132").Value = 1927.25
$ws.Range("K132").Value = 103601.319
$ws.Range("L132").Value = 5781.75
$ws.Range("M132").Value = -101071.319
$ws.Range("N132").Value = -10841.75
$ws.Range("H136").Value = 56988.633
$ws.Range("I136").Value = 3975.7778
$ws.Range("K136").Value = 11927.3334
$ws.Range("M136").Value = -9377.3334
